# Rename the sheet (Hoja1 -> "tiempo de ejecución")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "tiempo de ejecución"

# Update the two shared-string texts that changed.
# Set C25's comment first, then A19's header, so the shared-string pool
# ends up ordered the same way as in the target workbook.
$ws.Range("C25").Value = "(se utilizó n=16k como máxim n porque demoraba más de 8 minutos para 1 proceso)"
$ws.Range("A19").Value = "Tiempo de ejecución en segundos"

# Select A19:M29 and apply an all-around thin box border to the whole
# summary table (this also creates the now-visible empty bordered cells
# in rows 19, 20, 23-28 for columns B:M / A respectively).
$ws.Range("A19:M29").Borders.LineStyle = 1

# Re-create the selection/scroll state roughly matching the saved view.
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A19:M29").Select()

# Page setup: A4 paper, portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
